# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# D-column price strings that look numeric are written with a leading apostrophe
# so Excel keeps them as text (matching the sheet's inlineStr storage) instead of
# coercing them into floating-point numbers and dropping significant digits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.189.98'
$ws.Range("E2").Value = '  -1.27%  '

$ws.Range("D3").Value = '2.315.77'
$ws.Range("E3").Value = '  -2.29%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''312.91'
$ws.Range("E5").Value = '  -5.88%  '

$ws.Range("D6").Value = '''106.19'
$ws.Range("E6").Value = '  +5.21%  '

$ws.Range("E7").Value = '  -1.72%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").Value = '''0.610'
$ws.Range("E9").Value = '  -2.99%  '

$ws.Range("D10").Value = '''40.27'
$ws.Range("E10").Value = '  +0.72%  '

$ws.Range("D11").Value = '''0.0914'
$ws.Range("E11").Value = '  -1.08%  '

$ws.Range("D12").Value = '''8.28'
$ws.Range("E12").Value = '  -2.45%  '

$ws.Range("E13").Value = '  -0.18%  '

$ws.Range("D14").Value = '''0.981'
$ws.Range("E14").Value = '  -2.82%  '

$ws.Range("D15").Value = '''15.62'
$ws.Range("E15").Value = '  -5.22%  '

$ws.Range("D16").Value = '2.657.25'
$ws.Range("E16").Value = '  -2.59%  '

$ws.Range("D17").Value = '2.319.15'
$ws.Range("E17").Value = '  -2.11%  '

$ws.Range("D18").Value = '42.177.38'
$ws.Range("E18").Value = '  -1.34%  '

$ws.Range("E19").Value = '  -3.74%  '

$ws.Range("E20").Value = '  -1.75%  '

$ws.Range("D21").Value = '''74.68'
$ws.Range("E21").Value = '  -2.23%  '

$ws.Range("D22").Value = '''3.48'
$ws.Range("E22").Value = '  -8.23%  '

$ws.Range("D23").Value = '''256.32'
$ws.Range("E23").Value = '  -5.06%  '

$ws.Range("D24").Value = '''2.30'
$ws.Range("E24").Value = '  -0.78%  '

$ws.Range("D25").Value = '''9.29'
$ws.Range("E25").Value = '  -8.49%  '

$ws.Range("E26").Value = '  +0.39%  '

$ws.Range("D27").Value = '''11.01'
$ws.Range("E27").Value = '  -4.47%  '

$ws.Range("E28").Value = '  +3.24%  '

$ws.Range("D29").Value = '''22.81'
$ws.Range("E29").Value = '  -1.84%  '

$ws.Range("D30").Value = '''35.63'
$ws.Range("E30").Value = '  +0.72%  '

$ws.Range("D31").Value = '''0.0897'
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("D32").Value = '''163.12'
$ws.Range("E32").Value = '  -7.58%  '

$ws.Range("D33").Value = '''2.90'
$ws.Range("E33").Value = '  -6.13%  '

$ws.Range("D34").Value = '''5.86'
$ws.Range("E34").Value = '  -4.80%  '

# Row 35: Kaspa/Stellar swap places in the ranking
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '''0.119'
$ws.Range("E35").Value = '  +12.80%  '

# Row 36: Kaspa/Stellar swap places in the ranking
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = '''0.129'
$ws.Range("E36").Value = '  -2.64%  '

$ws.Range("D37").Value = '''4.54'
$ws.Range("E37").Value = '  -1.81%  '

$ws.Range("D38").Value = '''0.0353'
$ws.Range("E38").Value = '  -1.81%  '

$ws.Range("D39").Value = '''2.79'
$ws.Range("E39").Value = '  -6.34%  '

$ws.Range("D40").Value = '''3.65'
$ws.Range("E40").Value = '  -4.37%  '

$ws.Range("D41").Value = '''98.32'
$ws.Range("E41").Value = '  +7.41%  '

$ws.Range("D42").Value = '''1.47'
$ws.Range("E42").Value = '  -4.46%  '

$ws.Range("D43").Value = '''70.60'
$ws.Range("E43").Value = '  +0.55%  '

$ws.Range("D44").Value = '''0.230'
$ws.Range("E44").Value = '  -2.26%  '

$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("D46").Value = '''12.18'
$ws.Range("E46").Value = '  +2.26%  '

$ws.Range("D47").Value = '''111.73'
$ws.Range("E47").Value = '  -5.43%  '

$ws.Range("D48").Value = '''5.38'
$ws.Range("E48").Value = '  -2.42%  '

$ws.Range("D49").Value = '''9.00'
$ws.Range("E49").Value = '  -2.50%  '

$ws.Range("D50").Value = '''74.85'
$ws.Range("E50").Value = '  +5.95%  '

$ws.Range("D51").Value = '''1.27'
$ws.Range("E51").Value = '  -0.26%  '
